$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "CasesTab" row to "ParticipantsTab"
$ws.Range("A2").Value = "ParticipantsTab"

# Update the active selection to A2
$ws.Range("A2").Select()
